# Add two new BOM rows for heatsinks (pin header / PCB layout related BOM update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url5 = "https://www.mouser.de/ProductDetail/Advanced-Thermal-Solutions/ATS024024014-SF-8M?qs=9vOqFld9vZUhD4gLG7MQxw%3D%3D"
$url6 = "https://www.mouser.de/ProductDetail/Advanced-Thermal-Solutions/ATS024024010-SF-8I?qs=9vOqFld9vZWxOh9SEW0GcA%3D%3D"

# Fill the new cells in the same order the strings were first typed, so the
# shared-string table is built up in the original sequence.
$ws.Range("A5").Value = "heatsink 24x24x14mm"
$ws.Range("A6").Value = "heatsink 24x24x10mm"
$ws.Range("B5").Value = "ATS024024014-SF-8M"
$ws.Range("D5").Value = "Mouser"
$ws.Range("D6").Value = "Mouser"
$ws.Range("B6").Value = "ATS024024010-SF-8i"
$ws.Range("C5").Value = 5
$ws.Range("C6").Value = 5
$ws.Range("E5").Value = $url5
$ws.Range("E6").Value = $url6

# Turn E5/E6 into real hyperlinks, matching the style of existing link cells
$ws.Hyperlinks.Add($ws.Range("E5"), $url5) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E6"), $url6) | Out-Null

# Re-apply the existing hyperlink-cell format (centered "Hyperlink" style) from
# E1 so E5/E6 end up sharing the same style record as the other link cells.
$ws.Range("E1").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Move the active selection cursor as it ended up after the edit session
$ws.Range("C36").Select()
